$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current row 2 and row 3 values for the columns that need to be swapped
$d2 = $ws.Range("D2").Value()
$d3 = $ws.Range("D3").Value()

$n2 = $ws.Range("N2").Value()
$n3 = $ws.Range("N3").Value()

$o2 = $ws.Range("O2").Value()
$o3 = $ws.Range("O3").Value()

$p2 = $ws.Range("P2").Value()
$p3 = $ws.Range("P3").Value()

$q2 = $ws.Range("Q2").Value()
$q3 = $ws.Range("Q3").Value()

$r2 = $ws.Range("R2").Value()
$r3 = $ws.Range("R3").Value()

$s2 = $ws.Range("S2").Value()
$s3 = $ws.Range("S3").Value()

$t2 = $ws.Range("T2").Value()
$t3 = $ws.Range("T3").Value()

# Write back the swapped values: row2 gets row3's values, row3 gets row2's values
$ws.Range("D2").Value = $d3
$ws.Range("D3").Value = $d2

$ws.Range("N2").Value = $n3
$ws.Range("N3").Value = $n2

$ws.Range("O2").Value = $o3
$ws.Range("O3").Value = $o2

$ws.Range("P2").Value = $p3
$ws.Range("P3").Value = $p2

$ws.Range("Q2").Value = $q3
$ws.Range("Q3").Value = $q2

$ws.Range("R2").Value = $r3
$ws.Range("R3").Value = $r2

$ws.Range("S2").Value = $s3
$ws.Range("S3").Value = $s2

$ws.Range("T2").Value = $t3
$ws.Range("T3").Value = $t2
